$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before C (shifts old C:H -> D:I), inheriting formatting
# from the adjacent column as Excel normally does on a real column insert.
$ws.Columns("C:C").Insert()

# New header for the inserted column.
$ws.Range("C3").Value = "Accuracy after attack"

# New "raw accuracy after attack" values for the inserted column, row by row.
$ws.Range("C4").Value = 0.5
$ws.Range("C5").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 99.09999999999999
$ws.Range("C10").Value = 97.5
$ws.Range("C11").Value = 95.8
$ws.Range("C12").Value = 94
$ws.Range("C13").Value = 91.7
$ws.Range("C14").Value = 86.3
$ws.Range("C15").Value = 76.2
$ws.Range("C16").Value = 46.9
$ws.Range("C17").Value = 10
$ws.Range("C18").Value = 0.8
$ws.Range("C19").Value = 0.2
$ws.Range("C20").Value = 0
$ws.Range("C21").Value = 0
$ws.Range("C22").Value = 0
$ws.Range("C23").Value = 0
$ws.Range("C24").Value = 0
$ws.Range("C25").Value = 0
$ws.Range("C26").Value = 0.2
$ws.Range("C27").Value = 0.3
$ws.Range("C28").Value = 7.000000000000001
$ws.Range("C29").Value = 60.2
$ws.Range("C30").Value = 11.4
$ws.Range("C31").Value = 11
$ws.Range("C32").Value = 12.4
$ws.Range("C33").Value = 12.4
$ws.Range("C34").Value = 78.10000000000001
$ws.Range("C35").Value = 93.3
$ws.Range("C36").Value = 75.3
$ws.Range("C37").Value = 98.5
$ws.Range("C38").Value = 93.90000000000001
$ws.Range("C39").Value = 79.90000000000001
$ws.Range("C40").Value = 49.4
